$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "...desja gecter Mays..." -> "...desja gectee Mays..."
# The trailing "r" of "gecter" lives in its own run (rPr: rtl only,
# no color). Locate it precisely and swap just that character.
# ------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("desja gecter Mays", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rChar = $d.Range($find1.End - 5, $find1.End - 4)   # the "r" right before " Mays"
$rChar.Text = "e"

# ------------------------------------------------------------------
# Change 2: " Laquelle aplique fresche " -> " Laquelle apliquee fresche "
# Entirely inside a single run, plain replace is safe.
# ------------------------------------------------------------------
$d.Content.Find.Execute(" Laquelle aplique fresche ", $true, $false, $false, $false, $false, $true, 1, $false, " Laquelle apliquee fresche ", 2)

# ------------------------------------------------------------------
# Change 3: "quil est pre" -> "quil est pr" + new run "is"
# ("quil est pris"), where "is" must land in its OWN run carrying
# only <w:rtl val="0"/> (no explicit color), matching a sibling run
# elsewhere in the doc ("t" inside "...et animal..."). We borrow that
# run's formatting via FormattedText so the split run has no color
# override, then rewrite its text to "is".
# ------------------------------------------------------------------
$cleanSrc = $d.Content
$cleanSrc.Find.Execute("et animal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$cleanRun = $d.Range($cleanSrc.Start + 1, $cleanSrc.Start + 2)   # the "t" of "et"
$cleanFormat = $cleanRun.FormattedText

$find3 = $d.Content
$find3.Find.Execute("quil est pre", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailE = $d.Range($find3.End - 1, $find3.End)
$tailE.Text = ""
$insertPoint = $d.Range($find3.End - 1, $find3.End - 1)
$insertPoint.FormattedText = $cleanFormat
$newRun = $d.Range($find3.End - 1, $find3.End)
$newRun.Text = "is"

# ------------------------------------------------------------------
# Change 4: "est les meilleur" -> "est le meilleur"
# Entirely inside a single run, plain replace is safe.
# ------------------------------------------------------------------
$d.Content.Find.Execute("est les meilleur", $true, $false, $false, $false, $false, $true, 1, $false, "est le meilleur", 2)
